$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_HE_THONG")

# Update last_edited_time (column D) for rows 7-12 (Thang 7 .. Thang 2)
$ws.Range("D7").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D8").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D9").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D10").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D11").Value = "2024-07-18T15:58:00.000Z"
$ws.Range("D12").Value = "2024-07-18T15:58:00.000Z"

# Update row 7 (Thang 7) numeric property values due to the updated % format calc
$ws.Range("T7").Value = 37300000
$ws.Range("W7").Value = 216711000
$ws.Range("AA7").Value = 163107000
$ws.Range("AE7").Value = 379818000
$ws.Range("AH7").Value = 325818000
$ws.Range("AK7").Value = 54
$ws.Range("AN7").Value = 54000000
$ws.Range("AQ7").Value = 363118000
